$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting Code/Description/Definition to B/C/D
$ws.Range("A1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("A1").Value = "Version"

# Fill "1.0" for every data row (rows 2-8); write it as a formula that
# evaluates to the text "1.0", then paste-special values-only over itself
# so the stored cell is a plain text value (not a number, and without
# forcing a quote-prefix style change).
$ws.Range("A2:A8").Formula = '="1.0"'
$ws.Range("A2:A8").Copy()
$ws.Range("A2:A8").PasteSpecial(-4163)

# Match the updated sheetFormatPr baseColWidth attribute
$ws.StandardWidth = 10
